$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "0M"
$t.Cell(2,1).Range.Text = "0M"
$t.Cell(3,1).Range.Text = "0M"
$t.Cell(4,1).Range.Text = "412"

$t.Cell(6,1).Range.Text = "0.00068"
$t.Cell(7,1).Range.Text = "0.00023"
$t.Cell(8,1).Range.Text = "0.00005"
$t.Cell(9,1).Range.Text = "0.00044"
$t.Cell(10,1).Range.Text = "0.00046"
$t.Cell(11,1).Range.Text = "0.00048"
$t.Cell(12,1).Range.Text = "0.09465"

$t.Cell(44,1).Range.Text = "99.96"
$t.Cell(45,1).Range.Text = "0.09"
$t.Cell(46,1).Range.Text = "212"
